$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct the doses-administered figure for 2021-04-24 (row 125) ---
$ws.Range("B125").Value = 3580

# --- 2. Preserve the "latest day" (green/Good) formatting by moving it down
#        to the new last row (128) before we repaint row 127 ---
$ws.Range("A127:N127").Copy()
$ws.Range("A128:N128").PasteSpecial(-4122)   # xlPasteFormats

# Row 127 becomes a normal (yellow/Neutral) row, like row 126
$ws.Range("A126:N126").Copy()
$ws.Range("A127:N127").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 3. Correct the doses-administered figure for 2021-04-26 (row 127) ---
$ws.Range("B127").Value = 3582

# --- 4. Append the new day's data: 2021-04-27 (row 128) ---
$ws.Range("A128").Value = 44313
$ws.Range("B128").Value = 3537
$ws.Range("C128").Formula = "=(AVERAGE(B122:B128))"
$ws.Range("D128").Formula = "=D127-B128"
$ws.Range("E128").Formula = "=E127+B128"
$ws.Range("F128").Formula = "=(E128-G128)"
$ws.Range("G128").Value = 14306
$ws.Range("H128").Value = 28612
$ws.Range("I128").Formula = "=G128/2"
$ws.Range("J128").Value = 7153
$ws.Range("K128").Formula = "=D128/C128"
$ws.Range("L128").Formula = "=A128+K128"
$ws.Range("M128").Value = 476352
$ws.Range("N128").Value = 476352

# --- 5. Update view state to match: scrolled down one row, selection on J129 ---
$ws.Range("J129").Select() | Out-Null
$ActiveWindow = $excel.ActiveWindow
$ActiveWindow.ScrollRow = 115

$wb.Application.Calculate() | Out-Null
